$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old contents (but keep formatting) so the shared-strings table rebuilds cleanly
$ws.Range("A1:T13").ClearContents()

# Re-write header row (order of first appearance controls shared-string table order)
$ws.Cells.Item(1, 1).Value = "Sending cluster"
$ws.Cells.Item(1, 2).Value = "Ligand symbol"
$ws.Cells.Item(1, 3).Value = "Receptor symbol"
$ws.Cells.Item(1, 4).Value = "Target cluster"
$ws.Cells.Item(1, 5).Value = "Ligand-expressing cells"
$ws.Cells.Item(1, 6).Value = "Ligand detection rate"
$ws.Cells.Item(1, 7).Value = "Ligand average expression value"
$ws.Cells.Item(1, 8).Value = "Ligand total expression value"
$ws.Cells.Item(1, 9).Value = "Ligand derived specificity of average expression value"
$ws.Cells.Item(1, 10).Value = "Ligand derived specificity of total expression value"
$ws.Cells.Item(1, 11).Value = "Receptor-expressing cells"
$ws.Cells.Item(1, 12).Value = "Receptor detection rate"
$ws.Cells.Item(1, 13).Value = "Receptor average expression value"
$ws.Cells.Item(1, 14).Value = "Receptor total expression value"
$ws.Cells.Item(1, 15).Value = "Receptor derived specificity of average expression value"
$ws.Cells.Item(1, 16).Value = "Receptor derived specificity of total expression value"
$ws.Cells.Item(1, 17).Value = "Edge average expression weight"
$ws.Cells.Item(1, 18).Value = "Edge total expression weight"
$ws.Cells.Item(1, 19).Value = "Edge average expression derived specificity"
$ws.Cells.Item(1, 20).Value = "Edge total expression derived specificity"

# Re-write text columns A-D column-by-column so new shared strings are interned in the
# same order as the authoring tool (ECs, FAPs, M2, sCs, Sema5a, Plxnb3)
# Column A
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(17, 1).Value = "sCs"
# Column B
$ws.Cells.Item(2, 2).Value = "Sema5a"
$ws.Cells.Item(3, 2).Value = "Sema5a"
$ws.Cells.Item(4, 2).Value = "Sema5a"
$ws.Cells.Item(5, 2).Value = "Sema5a"
$ws.Cells.Item(6, 2).Value = "Sema5a"
$ws.Cells.Item(7, 2).Value = "Sema5a"
$ws.Cells.Item(8, 2).Value = "Sema5a"
$ws.Cells.Item(9, 2).Value = "Sema5a"
$ws.Cells.Item(10, 2).Value = "Sema5a"
$ws.Cells.Item(11, 2).Value = "Sema5a"
$ws.Cells.Item(12, 2).Value = "Sema5a"
$ws.Cells.Item(13, 2).Value = "Sema5a"
$ws.Cells.Item(14, 2).Value = "Sema5a"
$ws.Cells.Item(15, 2).Value = "Sema5a"
$ws.Cells.Item(16, 2).Value = "Sema5a"
$ws.Cells.Item(17, 2).Value = "Sema5a"
# Column C
$ws.Cells.Item(2, 3).Value = "Plxnb3"
$ws.Cells.Item(3, 3).Value = "Plxnb3"
$ws.Cells.Item(4, 3).Value = "Plxnb3"
$ws.Cells.Item(5, 3).Value = "Plxnb3"
$ws.Cells.Item(6, 3).Value = "Plxnb3"
$ws.Cells.Item(7, 3).Value = "Plxnb3"
$ws.Cells.Item(8, 3).Value = "Plxnb3"
$ws.Cells.Item(9, 3).Value = "Plxnb3"
$ws.Cells.Item(10, 3).Value = "Plxnb3"
$ws.Cells.Item(11, 3).Value = "Plxnb3"
$ws.Cells.Item(12, 3).Value = "Plxnb3"
$ws.Cells.Item(13, 3).Value = "Plxnb3"
$ws.Cells.Item(14, 3).Value = "Plxnb3"
$ws.Cells.Item(15, 3).Value = "Plxnb3"
$ws.Cells.Item(16, 3).Value = "Plxnb3"
$ws.Cells.Item(17, 3).Value = "Plxnb3"
# Column D
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(17, 4).Value = "sCs"

# Re-write numeric columns E-T
# Column E
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(17, 5).Value = 3
# Column F
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(17, 6).Value = 1
# Column G
$ws.Cells.Item(2, 7).Value = 1.143562333333333
$ws.Cells.Item(3, 7).Value = 1.143562333333333
$ws.Cells.Item(4, 7).Value = 1.143562333333333
$ws.Cells.Item(5, 7).Value = 1.143562333333333
$ws.Cells.Item(6, 7).Value = 40.163957
$ws.Cells.Item(7, 7).Value = 40.163957
$ws.Cells.Item(8, 7).Value = 40.163957
$ws.Cells.Item(9, 7).Value = 40.163957
$ws.Cells.Item(10, 7).Value = 0.002171666666666667
$ws.Cells.Item(11, 7).Value = 0.002171666666666667
$ws.Cells.Item(12, 7).Value = 0.002171666666666667
$ws.Cells.Item(13, 7).Value = 0.002171666666666667
$ws.Cells.Item(14, 7).Value = 5.751294333333334
$ws.Cells.Item(15, 7).Value = 5.751294333333334
$ws.Cells.Item(16, 7).Value = 5.751294333333334
$ws.Cells.Item(17, 7).Value = 5.751294333333334
# Column H
$ws.Cells.Item(2, 8).Value = 3.430687
$ws.Cells.Item(3, 8).Value = 3.430687
$ws.Cells.Item(4, 8).Value = 3.430687
$ws.Cells.Item(5, 8).Value = 3.430687
$ws.Cells.Item(6, 8).Value = 120.491871
$ws.Cells.Item(7, 8).Value = 120.491871
$ws.Cells.Item(8, 8).Value = 120.491871
$ws.Cells.Item(9, 8).Value = 120.491871
$ws.Cells.Item(10, 8).Value = 0.006515
$ws.Cells.Item(11, 8).Value = 0.006515
$ws.Cells.Item(12, 8).Value = 0.006515
$ws.Cells.Item(13, 8).Value = 0.006515
$ws.Cells.Item(14, 8).Value = 17.253883
$ws.Cells.Item(15, 8).Value = 17.253883
$ws.Cells.Item(16, 8).Value = 17.253883
$ws.Cells.Item(17, 8).Value = 17.253883
# Column I
$ws.Cells.Item(2, 9).Value = 0.02429958330097579
$ws.Cells.Item(3, 9).Value = 0.02429958330097579
$ws.Cells.Item(4, 9).Value = 0.02429958330097579
$ws.Cells.Item(5, 9).Value = 0.02429958330097579
$ws.Cells.Item(6, 9).Value = 0.8534448804146018
$ws.Cells.Item(7, 9).Value = 0.8534448804146018
$ws.Cells.Item(8, 9).Value = 0.8534448804146018
$ws.Cells.Item(9, 9).Value = 0.8534448804146018
$ws.Cells.Item(10, 9).Value = 0.00004614579680567107
$ws.Cells.Item(11, 9).Value = 0.00004614579680567107
$ws.Cells.Item(12, 9).Value = 0.00004614579680567107
$ws.Cells.Item(13, 9).Value = 0.00004614579680567107
$ws.Cells.Item(14, 9).Value = 0.1222093904876166
$ws.Cells.Item(15, 9).Value = 0.1222093904876166
$ws.Cells.Item(16, 9).Value = 0.1222093904876166
$ws.Cells.Item(17, 9).Value = 0.1222093904876166
# Column J
$ws.Cells.Item(2, 10).Value = 0.0242995833009758
$ws.Cells.Item(3, 10).Value = 0.0242995833009758
$ws.Cells.Item(4, 10).Value = 0.0242995833009758
$ws.Cells.Item(5, 10).Value = 0.0242995833009758
$ws.Cells.Item(6, 10).Value = 0.853444880414602
$ws.Cells.Item(7, 10).Value = 0.853444880414602
$ws.Cells.Item(8, 10).Value = 0.853444880414602
$ws.Cells.Item(9, 10).Value = 0.853444880414602
$ws.Cells.Item(10, 10).Value = 0.00004614579680567108
$ws.Cells.Item(11, 10).Value = 0.00004614579680567108
$ws.Cells.Item(12, 10).Value = 0.00004614579680567108
$ws.Cells.Item(13, 10).Value = 0.00004614579680567108
$ws.Cells.Item(14, 10).Value = 0.1222093904876167
$ws.Cells.Item(15, 10).Value = 0.1222093904876167
$ws.Cells.Item(16, 10).Value = 0.1222093904876167
$ws.Cells.Item(17, 10).Value = 0.1222093904876167
# Column K
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(17, 11).Value = 3
# Column L
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(17, 12).Value = 1
# Column M
$ws.Cells.Item(2, 13).Value = 0.367235
$ws.Cells.Item(3, 13).Value = 0.382402
$ws.Cells.Item(4, 13).Value = 1.055037666666667
$ws.Cells.Item(5, 13).Value = 1.515964666666667
$ws.Cells.Item(6, 13).Value = 0.367235
$ws.Cells.Item(7, 13).Value = 0.382402
$ws.Cells.Item(8, 13).Value = 1.055037666666667
$ws.Cells.Item(9, 13).Value = 1.515964666666667
$ws.Cells.Item(10, 13).Value = 0.367235
$ws.Cells.Item(11, 13).Value = 0.382402
$ws.Cells.Item(12, 13).Value = 1.055037666666667
$ws.Cells.Item(13, 13).Value = 1.515964666666667
$ws.Cells.Item(14, 13).Value = 0.367235
$ws.Cells.Item(15, 13).Value = 0.382402
$ws.Cells.Item(16, 13).Value = 1.055037666666667
$ws.Cells.Item(17, 13).Value = 1.515964666666667
# Column N
$ws.Cells.Item(2, 14).Value = 1.101705
$ws.Cells.Item(3, 14).Value = 1.147206
$ws.Cells.Item(4, 14).Value = 3.165113
$ws.Cells.Item(5, 14).Value = 4.547894
$ws.Cells.Item(6, 14).Value = 1.101705
$ws.Cells.Item(7, 14).Value = 1.147206
$ws.Cells.Item(8, 14).Value = 3.165113
$ws.Cells.Item(9, 14).Value = 4.547894
$ws.Cells.Item(10, 14).Value = 1.101705
$ws.Cells.Item(11, 14).Value = 1.147206
$ws.Cells.Item(12, 14).Value = 3.165113
$ws.Cells.Item(13, 14).Value = 4.547894
$ws.Cells.Item(14, 14).Value = 1.101705
$ws.Cells.Item(15, 14).Value = 1.147206
$ws.Cells.Item(16, 14).Value = 3.165113
$ws.Cells.Item(17, 14).Value = 4.547894
# Column O
$ws.Cells.Item(2, 15).Value = 0.1105916551411084
$ws.Cells.Item(3, 15).Value = 0.115159149071494
$ws.Cells.Item(4, 15).Value = 0.3177212460492045
$ws.Cells.Item(5, 15).Value = 0.4565279497381931
$ws.Cells.Item(6, 15).Value = 0.1105916551411084
$ws.Cells.Item(7, 15).Value = 0.115159149071494
$ws.Cells.Item(8, 15).Value = 0.3177212460492045
$ws.Cells.Item(9, 15).Value = 0.4565279497381931
$ws.Cells.Item(10, 15).Value = 0.1105916551411084
$ws.Cells.Item(11, 15).Value = 0.115159149071494
$ws.Cells.Item(12, 15).Value = 0.3177212460492045
$ws.Cells.Item(13, 15).Value = 0.4565279497381931
$ws.Cells.Item(14, 15).Value = 0.1105916551411084
$ws.Cells.Item(15, 15).Value = 0.115159149071494
$ws.Cells.Item(16, 15).Value = 0.3177212460492045
$ws.Cells.Item(17, 15).Value = 0.4565279497381931
# Column P
$ws.Cells.Item(2, 16).Value = 0.1105916551411084
$ws.Cells.Item(3, 16).Value = 0.115159149071494
$ws.Cells.Item(4, 16).Value = 0.3177212460492045
$ws.Cells.Item(5, 16).Value = 0.456527949738193
$ws.Cells.Item(6, 16).Value = 0.1105916551411084
$ws.Cells.Item(7, 16).Value = 0.115159149071494
$ws.Cells.Item(8, 16).Value = 0.3177212460492045
$ws.Cells.Item(9, 16).Value = 0.456527949738193
$ws.Cells.Item(10, 16).Value = 0.1105916551411084
$ws.Cells.Item(11, 16).Value = 0.115159149071494
$ws.Cells.Item(12, 16).Value = 0.3177212460492045
$ws.Cells.Item(13, 16).Value = 0.456527949738193
$ws.Cells.Item(14, 16).Value = 0.1105916551411084
$ws.Cells.Item(15, 16).Value = 0.115159149071494
$ws.Cells.Item(16, 16).Value = 0.3177212460492045
$ws.Cells.Item(17, 16).Value = 0.456527949738193
# Column Q
$ws.Cells.Item(2, 17).Value = 0.4199561134816667
$ws.Cells.Item(3, 17).Value = 0.4373005233913333
$ws.Cells.Item(4, 17).Value = 1.206501335847889
$ws.Cells.Item(5, 17).Value = 1.733600091464222
$ws.Cells.Item(6, 17).Value = 14.749610748895
$ws.Cells.Item(7, 17).Value = 15.358777484714
$ws.Cells.Item(8, 17).Value = 42.37448747738033
$ws.Cells.Item(9, 17).Value = 60.88713968551934
$ws.Cells.Item(10, 17).Value = 0.0007975120083333334
$ws.Cells.Item(11, 17).Value = 0.0008304496766666666
$ws.Cells.Item(12, 17).Value = 0.002291190132777778
$ws.Cells.Item(13, 17).Value = 0.003292169934444445
$ws.Cells.Item(14, 17).Value = 2.112076574501667
$ws.Cells.Item(15, 17).Value = 2.199306455655333
$ws.Cells.Item(16, 17).Value = 6.067832153753222
$ws.Cells.Item(17, 17).Value = 8.718758996933557
# Column R
$ws.Cells.Item(2, 18).Value = 3.779605021335
$ws.Cells.Item(3, 18).Value = 3.935704710522
$ws.Cells.Item(4, 18).Value = 10.858512022631
$ws.Cells.Item(5, 18).Value = 15.602400823178
$ws.Cells.Item(6, 18).Value = 132.746496740055
$ws.Cells.Item(7, 18).Value = 138.228997362426
$ws.Cells.Item(8, 18).Value = 381.370387296423
$ws.Cells.Item(9, 18).Value = 547.984257169674
$ws.Cells.Item(10, 18).Value = 0.007177608075
$ws.Cells.Item(11, 18).Value = 0.00747404709
$ws.Cells.Item(12, 18).Value = 0.020620711195
$ws.Cells.Item(13, 18).Value = 0.02962952941
$ws.Cells.Item(14, 18).Value = 19.008689170515
$ws.Cells.Item(15, 18).Value = 19.793758100898
$ws.Cells.Item(16, 18).Value = 54.610489383779
$ws.Cells.Item(17, 18).Value = 78.46883097240202
# Column S
$ws.Cells.Item(2, 19).Value = 0.00268733113649415
$ws.Cells.Item(3, 19).Value = 0.002798319335732258
$ws.Cells.Item(4, 19).Value = 0.007720493884862471
$ws.Cells.Item(5, 19).Value = 0.01109343894388691
$ws.Cells.Item(6, 19).Value = 0.0943838818967561
$ws.Cells.Item(7, 19).Value = 0.09828198620796853
$ws.Cells.Item(8, 19).Value = 0.2711575708396416
$ws.Cells.Item(9, 19).Value = 0.3896214414702355
$ws.Cells.Item(10, 19).Value = 0.000005103340046544435
$ws.Cells.Item(11, 19).Value = 0.000005314110693367149
$ws.Cells.Item(12, 19).Value = 0.00001466150006103122
$ws.Cells.Item(13, 19).Value = 0.00002106684600472827
$ws.Cells.Item(14, 19).Value = 0.01351533876781155
$ws.Cells.Item(15, 19).Value = 0.01407352941709987
$ws.Cells.Item(16, 19).Value = 0.03882851982463936
$ws.Cells.Item(17, 19).Value = 0.05579200247806586
# Column T
$ws.Cells.Item(2, 20).Value = 0.002687331136494151
$ws.Cells.Item(3, 20).Value = 0.002798319335732259
$ws.Cells.Item(4, 20).Value = 0.007720493884862473
$ws.Cells.Item(5, 20).Value = 0.01109343894388691
$ws.Cells.Item(6, 20).Value = 0.09438388189675613
$ws.Cells.Item(7, 20).Value = 0.09828198620796856
$ws.Cells.Item(8, 20).Value = 0.2711575708396417
$ws.Cells.Item(9, 20).Value = 0.3896214414702355
$ws.Cells.Item(10, 20).Value = 0.000005103340046544435
$ws.Cells.Item(11, 20).Value = 0.00000531411069336715
$ws.Cells.Item(12, 20).Value = 0.00001466150006103122
$ws.Cells.Item(13, 20).Value = 0.00002106684600472827
$ws.Cells.Item(14, 20).Value = 0.01351533876781155
$ws.Cells.Item(15, 20).Value = 0.01407352941709988
$ws.Cells.Item(16, 20).Value = 0.03882851982463938
$ws.Cells.Item(17, 20).Value = 0.05579200247806586
